$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.63867
$ws.Range("H2").Value = 1.91601
$ws.Range("I2").Value = 0.02162993170452444
$ws.Range("J2").Value = 0.02162993170452444
$ws.Range("M2").Value = 0.491263
$ws.Range("N2").Value = 1.473789
$ws.Range("O2").Value = 0.08507192955174298
$ws.Range("P2").Value = 0.08507192955174299
$ws.Range("Q2").Value = 0.31375494021
$ws.Range("R2").Value = 2.82379446189
$ws.Range("S2").Value = 0.001840100026176315
$ws.Range("T2").Value = 0.001840100026176315

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.63867
$ws.Range("H3").Value = 1.91601
$ws.Range("I3").Value = 0.02162993170452444
$ws.Range("J3").Value = 0.02162993170452444
$ws.Range("M3").Value = 4.669265333333333
$ws.Range("N3").Value = 14.007796
$ws.Range("O3").Value = 0.8085758778815603
$ws.Range("P3").Value = 0.8085758778815605
$ws.Range("Q3").Value = 2.982119690439999
$ws.Range("R3").Value = 26.83907721396
$ws.Range("S3").Value = 0.01748944101650404
$ws.Range("T3").Value = 0.01748944101650405

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.63867
$ws.Range("H4").Value = 1.91601
$ws.Range("I4").Value = 0.02162993170452444
$ws.Range("J4").Value = 0.02162993170452444
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6141496666666667
$ws.Range("N4").Value = 1.842449
$ws.Range("O4").Value = 0.1063521925666967
$ws.Range("P4").Value = 0.1063521925666967
$ws.Range("Q4").Value = 0.39223896761
$ws.Range("R4").Value = 3.53015070849
$ws.Range("S4").Value = 0.00230039066184408
$ws.Range("T4").Value = 0.002300390661844081

# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("I5").Value = 0.9490095874171892
$ws.Range("J5").Value = 0.9490095874171893
$ws.Range("M5").Value = 0.491263
$ws.Range("N5").Value = 1.473789
$ws.Range("O5").Value = 0.08507192955174298
$ws.Range("P5").Value = 0.08507192955174299
$ws.Range("Q5").Value = 13.765944822494
$ws.Range("R5").Value = 123.893503402446
$ws.Range("S5").Value = 0.08073407676468379
$ws.Range("T5").Value = 0.0807340767646838

# Row 6
$ws.Range("D6").Value = "FAPs"
$ws.Range("I6").Value = 0.9490095874171892
$ws.Range("J6").Value = 0.9490095874171893
$ws.Range("M6").Value = 4.669265333333333
$ws.Range("N6").Value = 14.007796
$ws.Range("O6").Value = 0.8085758778815603
$ws.Range("P6").Value = 0.8085758778815605
$ws.Range("Q6").Value = 130.8399959700826
$ws.Range("R6").Value = 1177.559963730744
$ws.Range("S6").Value = 0.7673462602638711
$ws.Range("T6").Value = 0.7673462602638714

# Row 7
$ws.Range("D7").Value = "MuSCs"
$ws.Range("I7").Value = 0.9490095874171892
$ws.Range("J7").Value = 0.9490095874171893
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6141496666666667
$ws.Range("N7").Value = 1.842449
$ws.Range("O7").Value = 0.1063521925666967
$ws.Range("P7").Value = 0.1063521925666967
$ws.Range("Q7").Value = 17.20941822218733
$ws.Range("R7").Value = 154.884763999686
$ws.Range("S7").Value = 0.1009292503886343
$ws.Range("T7").Value = 0.1009292503886343

# Row 8
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 0.866931
$ws.Range("H8").Value = 2.600793
$ws.Range("I8").Value = 0.02936048087828625
$ws.Range("J8").Value = 0.02936048087828625
$ws.Range("M8").Value = 0.491263
$ws.Range("N8").Value = 1.473789
$ws.Range("O8").Value = 0.08507192955174298
$ws.Range("P8").Value = 0.08507192955174299
$ws.Range("Q8").Value = 0.425891123853
$ws.Range("R8").Value = 3.833020114677
$ws.Range("S8").Value = 0.002497752760882864
$ws.Range("T8").Value = 0.002497752760882865

# Row 9
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 0.866931
$ws.Range("H9").Value = 2.600793
$ws.Range("I9").Value = 0.02936048087828625
$ws.Range("J9").Value = 0.02936048087828625
$ws.Range("M9").Value = 4.669265333333333
$ws.Range("N9").Value = 14.007796
$ws.Range("O9").Value = 0.8085758778815603
$ws.Range("P9").Value = 0.8085758778815605
$ws.Range("Q9").Value = 4.047930864692
$ws.Range("R9").Value = 36.431377782228
$ws.Range("S9").Value = 0.02374017660118507
$ws.Range("T9").Value = 0.02374017660118507

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 0.866931
$ws.Range("H10").Value = 2.600793
$ws.Range("I10").Value = 0.02936048087828625
$ws.Range("J10").Value = 0.02936048087828625
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6141496666666667
$ws.Range("N10").Value = 1.842449
$ws.Range("O10").Value = 0.1063521925666967
$ws.Range("P10").Value = 0.1063521925666967
$ws.Range("Q10").Value = 0.532425384673
$ws.Range("R10").Value = 4.791828462057
$ws.Range("S10").Value = 0.003122551516218314
$ws.Range("T10").Value = 0.003122551516218314
